# Applies the "Add files via upload" edit to the ModeloEntidadRelacionDB
# workbook:
#   - E3 header renamed: id_persona_owner -> id_usuario_owner
#   - The "Usuario" table (rows 33-37) gains two new columns:
#       * a new "id_persona" column right after "id"
#       * a new "password" column right before "email"
#     and the table's merged title cell grows from B33:G33 to B33:I33
#   - The hyperlink on the email cell follows the email column (D35 -> F35)
#   - The "Log" table's "id_persona" column header is renamed to "id_user"
#   - Selection moves to E3 (the cell that was actually edited)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shifts the block [fromCol..toCol] in $row one column to the right,
# carrying values AND formatting (Cut preserves both). Cells are moved
# starting from the rightmost column so nothing gets clobbered before it
# is relocated. Leaves $fromCol blank/empty in $row.
function Shift-RowRight {
    param($row, $fromCol, $toCol)
    for ($c = $toCol; $c -ge $fromCol; $c--) {
        $src = $ws.Cells.Item($row, $c)
        $dst = $ws.Cells.Item($row, $c + 1)
        $src.Cut($dst) | Out-Null
    }
}

# --- 1. Rename the "Vector" table header cell E3 -------------------------
$ws.Range("E3").Value = "id_usuario_owner"

# --- 2. Make room for the two new "Usuario" table columns -----------------
# Old layout (row 34 header): B=id C=Usuario D=email E=id_usuario_padre
#                              F=idRolUsuario G=fechaRegistro
# New layout:                 B=id C=id_persona D=Usuario E=password
#                              F=email G=id_usuario_padre H=idRolUsuario
#                              I=fechaRegistro
foreach ($r in 33..37) {
    # Free column C: shift C..G right to D..H
    Shift-RowRight $r 3 7
    # Free column E (now holding the old "Usuario"/data value, which
    # should stay put) for the new "password" column: shift E..H right to F..I
    Shift-RowRight $r 5 8
}

# --- 3. Grow the merged title cell from B33:G33 to B33:I33 ----------------
$ws.Range("B33:G33").UnMerge() | Out-Null
$ws.Range("B33:I33").Merge() | Out-Null

# --- 4. Re-point the hyperlink that used to sit on the email cell --------
# The email text/format already moved from D35 to F35 via the cut-shift
# above, but the hyperlink object itself stays registered on the old
# address, so drop it and recreate it on F35.
foreach ($h in @($ws.Hyperlinks)) {
    $h.Delete()
}
$ws.Hyperlinks.Add($ws.Range("F35"), "mailto:wisrovi.rodriguez@gmail.com", "", "", "wisrovi.rodriguez@gmail.com") | Out-Null

# --- 5. Fill in the new header + data cells -------------------------------
$ws.Range("C34").Value = "id_persona"
$ws.Range("E34").Value = "password"

$ws.Range("C35").Value = 1
$ws.Range("E35").Value = "awefa"

# --- 6. Rename the "Log" table's id_persona header to id_user ------------
$ws.Range("C62").Value = "id_user"

# --- 7. Move the active selection to the cell that was edited ------------
$ws.Range("E3").Select() | Out-Null
